# break out stock.yaml completed
# Applies the "chartink_screener" breakout sheet update:
#   - E86/E87/E88 (bsecode) converted from text to true numeric values
#   - two new trailing rows (89, 90) appended with the 11:35:37 screener run

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10per change")

# --- Fix up existing rows 86-88: bsecode column (E) becomes numeric ---
$ws.Cells.Item(86, 5).Value = 533274
$ws.Cells.Item(87, 5).Value = 533273
$ws.Cells.Item(88, 5).Value = 526371

# --- Row 89: PRESTIGE ---
$ws.Cells.Item(89, 1).Value = "27/06/2024 11:35:37"
$ws.Cells.Item(89, 2).Value = 1
$ws.Cells.Item(89, 3).Value = "PRESTIGE"
$ws.Cells.Item(89, 4).Value = "Prestige Estates Projects Limited"
# bsecode stays text here (unlike rows 86-88) - force text, then drop the
# number-format style trace so the cell ends up plain/unstyled like the rest
$ws.Cells.Item(89, 5).NumberFormat = "@"
$ws.Cells.Item(89, 5).Value = "533274"
$ws.Cells.Item(89, 5).Style = "Normal"
$ws.Cells.Item(89, 6).Value = -3.73
$ws.Cells.Item(89, 7).Value = 1855.1
$ws.Cells.Item(89, 8).Value = 1399189

# --- Row 90: NMDC ---
$ws.Cells.Item(90, 1).Value = "27/06/2024 11:35:37"
$ws.Cells.Item(90, 2).Value = 2
$ws.Cells.Item(90, 3).Value = "NMDC"
$ws.Cells.Item(90, 4).Value = "Nmdc Limited"
$ws.Cells.Item(90, 5).NumberFormat = "@"
$ws.Cells.Item(90, 5).Value = "526371"
$ws.Cells.Item(90, 5).Style = "Normal"
$ws.Cells.Item(90, 6).Value = -1.5
$ws.Cells.Item(90, 7).Value = 245.5
$ws.Cells.Item(90, 8).Value = 20047167
